$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-83
# from serial date 45221 (2023-10-22) to 45224 (2023-10-25)
for ($r = 2; $r -le 83; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45224
}
